$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Demo Fund 2"
$ws.Range("B2").Value = "TSTF2 Port Co 3"
$ws.Range("D2").Value = 10000000

$ws.Range("C2").Select()
